# "Template info now in Excel"
# Rework the credential placeholder rows (3-7) on the Summary sheet so the
# labels and values describe the actual template fields (Username, Password,
# Payer, Billing Provider) instead of the old Insurance/username/password
# placeholder rows, and widen column B to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# The sheet is protected - unprotect so the cells can be edited, then restore
# protection at the end.
$ws.Unprotect()

# Row 3: was Insurance / Insurance -> Username / Username
$ws.Range("A3").Value = "Username"
$ws.Range("B3").Value = "Username"

# Row 4: label stays "Password"; value was the lowercase "username" hint -> "Password"
$ws.Range("A4").Value = "Password"
$ws.Range("B4").Value = "Password"

# Row 5: was Password / password -> Payer / Payer
$ws.Range("A5").Value = "Payer"
$ws.Range("B5").Value = "Payer"

# Row 6: was Payer / (blank) -> Billing Provider / Billing Provider
$ws.Range("A6").Value = "Billing Provider"
$ws.Range("B6").Value = "Billing Provider"

# Row 7: was Billing Provider / (blank) -> now an empty spacer row, matching
# the formatting already used for the blank divider rows (row 8 / row 19).
$ws.Range("A8").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("A7:B7").ClearContents()

# Widen column B so the new, longer labels/values fit (was 21.08984375, auto
# best-fit; now an explicit width of roughly 30.63).
$ws.Columns.Item(2).ColumnWidth = 29.83

# Restore sheet protection.
$ws.Protect()
